$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# New header for column F
$ws.Range("F1").Value = "dV/dt"

# New column F: dV/dt = C*A*1000/10^6, rows 2..25
for ($r = 2; $r -le 25; $r++) {
    $ws.Range("F$r").Formula = "=C$r*A$r*1000/10^6"
    $ws.Range("F$r").NumberFormat = "0.00E+00"
}

# Update the view's selection / scroll position to match the edited area
$ws.Range("F11").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 5
$win.ScrollColumn = 1
